# Generate Report for Handoff
#
# The localization status report is regenerated: the "Status" columns move
# from "Handed back: in sync with en-US" to "Ready for handoff", and the
# associated timestamps advance a little.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps bump forward ---
$wsOverview.Range("G2").Value = "2016-08-15 09:13:29"
$wsZhCn.Range("H2").Value = "2016-08-15 09:13:24"
$wsDeDe.Range("H2").Value = "2016-08-15 09:13:29"

# --- Re-fit the now-shorter Status columns to their new content ---
# "Ready for handoff" is noticeably shorter than "Handed back: in sync with
# en-US", so the Status column on every sheet narrows to match (same target
# width everywhere, since it's driven by the same status text). AutoFit first
# (so the columns track the new text like a real handoff report would), then
# pin the resulting width explicitly.
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZhCn.Columns.Item(3).AutoFit() | Out-Null
$wsDeDe.Columns.Item(3).AutoFit() | Out-Null

$newStatusColWidth = 16.3
$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth

